$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous content entirely so stale F/G columns and old rows go away
$ws.Cells.Clear()

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 10917.3,              11579, 9364,  0.2509726285934448),
    @(1, 11066.46666666667,    11859, 9396,  0.2598573048909505),
    @(2, 10921.36666666667,    11558, 9991,  0.3281592528025309),
    @(3, 11558.53333333333,    12256, 10469, 0.3025633017222086),
    @(4, 10512.93333333333,    11319, 9611,  0.3081936359405518),
    @(5, 11466.9,              12369, 9922,  0.298317281405131),
    @(6, 11013.23333333333,    11835, 10138, 0.3340381781260173),
    @(7, 10743.9,              11548, 10143, 0.3538715442021688),
    @(8, 10653.63333333333,    11415, 9151,  0.2255171696345011),
    @(9, 10609.66666666667,    11597, 9665,  0.2707310358683268)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
